# Weekly update: two new daily price records for "Espárragos" arrived for
# Vega Modelo de Temuco, so they are inserted at the top of the data block
# (right after the header row) and the rest of the rows shift down by two.
#
# Net effect vs. the previous sheet:
#   - old row 10 .. old row 16  ->  new row 12 .. new row 18 (unchanged)
#   - two brand-new rows are inserted as the new row 10 and row 11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data rows (old 10-16) down by inserting two new rows
# above row 10 (just like pasting two new daily readings at the top of the
# list, above the most recent previously-recorded rows).
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# --- New row 10 ---------------------------------------------------------
$ws.Cells.Item(10, 1).Value  = 10
$ws.Cells.Item(10, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value  = "La Araucanía"
$ws.Cells.Item(10, 4).Value  = 44466
$ws.Cells.Item(10, 5).Value  = 9
$ws.Cells.Item(10, 6).Value  = 300000000
$ws.Cells.Item(10, 7).Value  = "Espárragos"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 300
$ws.Cells.Item(10, 11).Value = 2000
$ws.Cells.Item(10, 12).Value = 2000
$ws.Cells.Item(10, 13).Value = 2000
$ws.Cells.Item(10, 14).Value = "$/kilo"
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 2000
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# --- New row 11 ---------------------------------------------------------
$ws.Cells.Item(11, 1).Value  = 10
$ws.Cells.Item(11, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value  = "La Araucanía"
$ws.Cells.Item(11, 4).Value  = 44466
$ws.Cells.Item(11, 5).Value  = 9
$ws.Cells.Item(11, 6).Value  = 300000000
$ws.Cells.Item(11, 7).Value  = "Espárragos"
$ws.Cells.Item(11, 8).Value  = "Sin especificar"
$ws.Cells.Item(11, 9).Value  = "Segunda"
$ws.Cells.Item(11, 10).Value = 50
$ws.Cells.Item(11, 11).Value = 1500
$ws.Cells.Item(11, 12).Value = 1500
$ws.Cells.Item(11, 13).Value = 1500
$ws.Cells.Item(11, 14).Value = "$/kilo"
$ws.Cells.Item(11, 15).Value = "Región del Maule"
$ws.Cells.Item(11, 16).Value = 1500
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = "Hortaliza"
